# Auto-generated edit script applying the diff changes
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)
$ws4 = $wb.Worksheets.Item(4)

$ws1.Range("F4").Value = 8494
$ws1.Range("F5").Value = 8494
$ws1.Range("F6").Value = 559
$ws1.Range("F7").Value = 7506
$ws1.Range("F9").Value = 628
$ws1.Range("F10").Value = 522
$ws1.Range("C11").Value = '北京·排球少年同好嘉年华4th'
$ws1.Range("D11").Value = '石景山路68号 北京首钢会展中心'
$ws1.Range("E11").Value = '2024.07.20 09:00-07.21 17:00'
$ws1.Range("F11").Value = 763
$ws1.Range("G11").Value = 90
$ws1.Range("H11").Value = 'https://show.bilibili.com/platform/detail.html?id=85947'
$ws1.Range("I11").Value = '//i2.hdslb.com/bfs/openplatform/202405/bTK0cxNF1716260812082.jpeg'
$ws1.Range("F14").Value = 126
$ws1.Range("F15").Value = 178
$ws1.Range("F16").Value = 12454
$ws1.Range("F18").Value = 20
$ws1.Range("F19").Value = 2623
$ws1.Range("F20").Value = 3787
$ws1.Range("F21").Value = 58
$ws1.Range("F22").Value = 47
$ws1.Range("F23").Value = 3038
$ws1.Range("F24").Value = 11
$ws1.Range("F25").Value = 130
$ws1.Range("F27").Value = 21
$ws1.Range("F28").Value = 42
$ws1.Range("F29").Value = 3409
$ws1.Range("F31").Value = 349
$ws1.Range("F32").Value = 1766
$ws1.Range("F34").Value = 151
$ws1.Range("F35").Value = 6165
$ws1.Range("F36").Value = 8
$ws1.Range("F37").Value = 108
$ws1.Range("F38").Value = 1878
$ws1.Range("F40").Value = 48
$ws1.Range("F41").Value = 940
$ws1.Range("F43").Value = 183
$ws1.Range("F46").Value = 1127
$ws1.Range("F47").Value = 1117
$ws1.Range("F48").Value = 1628
$ws1.Range("F49").Value = 32
$ws1.Range("F50").Value = 125
$ws2.Range("G7").Value = '不可售'
$ws2.Range("F16").Value = 112
$ws2.Range("F22").Value = 80
$ws3.Range("F2").Value = 358
$ws3.Range("F3").Value = 509
$ws3.Range("F4").Value = 18
$ws4.Range("F6").Value = 358
$ws4.Range("F7").Value = 509
$ws4.Range("F9").Value = 8494
$ws4.Range("C10").Value = '北京·AINI二次元派对【免票展会】'
$ws4.Range("D10").Value = '天竺镇裕翔路99号 北京欧陆时尚购物中心'
$ws4.Range("E10").Value = '2024.07.20 16:00-07.21 19:00'
$ws4.Range("F10").Value = 559
$ws4.Range("G10").Value = 50
$ws4.Range("H10").Value = 'https://show.bilibili.com/platform/detail.html?id=89134'
$ws4.Range("I10").Value = '//i2.hdslb.com/bfs/openplatform/202407/lb2k8yDD1720678733848.jpeg'
$ws4.Range("C11").Value = '北京·IDO动漫游戏嘉年华46th'
$ws4.Range("D11").Value = '京沈路与天北路交汇处西北角 中国国际展览中心新馆'
$ws4.Range("E11").Value = '2024.07.20 09:30-07.21 17:00'
$ws4.Range("F11").Value = 7506
$ws4.Range("G11").Value = 95
$ws4.Range("H11").Value = 'https://show.bilibili.com/platform/detail.html?id=83716'
$ws4.Range("I11").Value = '//i2.hdslb.com/bfs/openplatform/202405/9CAdQvG71716812495452.jpeg'
$ws4.Range("F12").Value = 7506
$ws4.Range("C13").Value = '北京·原神同人嘉年华10th'
$ws4.Range("D13").Value = '石景山路68号 北京首钢会展中心'
$ws4.Range("E13").Value = '2024.07.20 09:00-07.21 17:00'
$ws4.Range("F13").Value = 628
$ws4.Range("G13").Value = 90
$ws4.Range("H13").Value = 'https://show.bilibili.com/platform/detail.html?id=86012'
$ws4.Range("I13").Value = '//i1.hdslb.com/bfs/openplatform/202405/ulMhJXc61716260154833.jpeg'
$ws4.Range("C14").Value = '北京·国乙同好嘉年华9th'
$ws4.Range("F14").Value = 522
$ws4.Range("H14").Value = 'https://show.bilibili.com/platform/detail.html?id=86011'
$ws4.Range("I14").Value = '//i0.hdslb.com/bfs/openplatform/202405/AB4NkZsl1716272246698.jpeg'
$ws4.Range("B15").Value = '2024-07-21'
$ws4.Range("C15").Value = '【大会员提前抢】北京·ICOS内场-青柳尊哉'
$ws4.Range("E15").Value = '2024.07.21 09:00-07.21 17:00'
$ws4.Range("F15").Value = 237
$ws4.Range("G15").Value = 598
$ws4.Range("H15").Value = 'https://show.bilibili.com/platform/detail.html?id=86904'
$ws4.Range("I15").Value = '//i2.hdslb.com/bfs/openplatform/202406/WgYg2oTS1717576349209.jpeg'
$ws4.Range("C16").Value = '北京·Summer Overture'
$ws4.Range("D16").Value = '朝阳北路甲27号菁英梦谷·常营文创产业园南门B5座 WeShow Live 北京'
$ws4.Range("E16").Value = '2024.07.21 12:00-07.21 19:00'
$ws4.Range("F16").Value = 261
$ws4.Range("G16").Value = 78
$ws4.Range("H16").Value = 'https://show.bilibili.com/platform/detail.html?id=87481'
$ws4.Range("I16").Value = '//i1.hdslb.com/bfs/openplatform/202406/dP7KKEIk1718608495643.png'
$ws4.Range("C17").Value = '北京·航海王（ONE PIECE）25周年巡展'
$ws4.Range("D17").Value = '酒仙桥路2号北京798艺术区A区 北京798艺术区'
$ws4.Range("E17").Value = '2024.07.21 10:00-10.27 19:00'
$ws4.Range("F17").Value = 178
$ws4.Range("G17").Value = 98
$ws4.Range("H17").Value = 'https://show.bilibili.com/platform/detail.html?id=89233'
$ws4.Range("I17").Value = '//i2.hdslb.com/bfs/openplatform/202407/WxL0mO9g1721011505489.png'
$ws4.Range("B18").Value = '2024-07-23'
$ws4.Range("C18").Value = '北京·巴西浪漫风情——手风琴大满贯音乐家道格拉斯·博尔萨蒂专场音乐会'
$ws4.Range("D18").Value = '复兴路69号号3号楼6层601 爱乐汇艺术空间(五棵松万达广场)'
$ws4.Range("E18").Value = '2024.07.23 19:30-07.23 21:00'
$ws4.Range("F18").Value = 4
$ws4.Range("G18").Value = 140
$ws4.Range("H18").Value = 'https://show.bilibili.com/platform/detail.html?id=86922'
$ws4.Range("I18").Value = '//i1.hdslb.com/bfs/openplatform/202405/i14RABlz1716527544509.jpeg'
$ws4.Range("B19").Value = '2024-07-26'
$ws4.Range("C19").Value = '北京·第17届IJOY漫展xCGF游戏节'
$ws4.Range("D19").Value = '天辰东路7号 北京国家会议中心'
$ws4.Range("E19").Value = '2024.07.26 09:00-07.28 17:00'
$ws4.Range("F19").Value = 12455
$ws4.Range("G19").Value = 85
$ws4.Range("H19").Value = 'https://show.bilibili.com/platform/detail.html?id=84088'
$ws4.Range("I19").Value = '//i0.hdslb.com/bfs/openplatform/202404/EiPIQJ7R1712817059082.jpeg'
$ws4.Range("B20").Value = '2024-07-27'
$ws4.Range("C20").Value = '北京·第17届IJOY漫展【新田惠海专场见面会】'
$ws4.Range("E20").Value = '2024.07.27 11:00-07.27 14:10'
$ws4.Range("F20").Value = 108
$ws4.Range("G20").Value = 628
$ws4.Range("H20").Value = 'https://show.bilibili.com/platform/detail.html?id=87626'
$ws4.Range("I20").Value = '//i1.hdslb.com/bfs/openplatform/202406/NIBORFE21718696357203.jpeg'
$ws4.Range("F22").Value = 20
$ws4.Range("F23").Value = 2623
$ws4.Range("F24").Value = 2623
$ws4.Range("F25").Value = 3787
$ws4.Range("F26").Value = 11
$ws4.Range("F27").Value = 130
$ws4.Range("F29").Value = 21
$ws4.Range("F30").Value = 42
$ws4.Range("F32").Value = 3409
$ws4.Range("F33").Value = 349
$ws4.Range("F34").Value = 1766
$ws4.Range("F36").Value = 151
$ws4.Range("F37").Value = 6165
$ws4.Range("F38").Value = 80
$ws4.Range("F39").Value = 8
$ws4.Range("F40").Value = 108
$ws4.Range("F41").Value = 1878
$ws4.Range("F44").Value = 48
$ws4.Range("F45").Value = 940
$ws4.Range("F46").Value = 183
$ws4.Range("F48").Value = 1127
$ws4.Range("F49").Value = 1117
$ws4.Range("F50").Value = 1628
$ws4.Range("F51").Value = 32
$ws4.Range("F52").Value = 125
